$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 95, pushing the existing rows 95-109 down to 97-111.
$ws.Rows("95:96").Insert()

# Fill in the new row 95 with the new weekly record ($/caja 20 kilos, Primera).
$ws.Cells.Item(95, 1).Value = 1
$ws.Cells.Item(95, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(95, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(95, 4).Value = 44736
$ws.Cells.Item(95, 5).Value = 15
$ws.Cells.Item(95, 6).Value = 100112036
$ws.Cells.Item(95, 7).Value = "Caigua"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 120
$ws.Cells.Item(95, 11).Value = 8000
$ws.Cells.Item(95, 12).Value = 9000
$ws.Cells.Item(95, 13).Value = 8500
$ws.Cells.Item(95, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(95, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(95, 16).Value = 425
$ws.Cells.Item(95, 17).Value = 20
$ws.Cells.Item(95, 18).Value = "Hortaliza"

# Fill in the new row 96 with the new weekly record ($/caja 20 kilos, Segunda).
$ws.Cells.Item(96, 1).Value = 1
$ws.Cells.Item(96, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(96, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(96, 4).Value = 44736
$ws.Cells.Item(96, 5).Value = 15
$ws.Cells.Item(96, 6).Value = 100112036
$ws.Cells.Item(96, 7).Value = "Caigua"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Segunda"
$ws.Cells.Item(96, 10).Value = 140
$ws.Cells.Item(96, 11).Value = 6000
$ws.Cells.Item(96, 12).Value = 7000
$ws.Cells.Item(96, 13).Value = 6500
$ws.Cells.Item(96, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(96, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(96, 16).Value = 325
$ws.Cells.Item(96, 17).Value = 20
$ws.Cells.Item(96, 18).Value = "Hortaliza"

# Add two new rows (110 and 111) at the end, duplicating the old last two rows
# (which are now at 108-109 after the shift above).
$ws.Cells.Item(110, 1).Value = 1
$ws.Cells.Item(110, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(110, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(110, 4).Value = 44442
$ws.Cells.Item(110, 4).NumberFormat = $ws.Cells.Item(109, 4).NumberFormat
$ws.Cells.Item(110, 5).Value = 15
$ws.Cells.Item(110, 6).Value = 100112036
$ws.Cells.Item(110, 7).Value = "Caigua"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 140
$ws.Cells.Item(110, 11).Value = 9000
$ws.Cells.Item(110, 12).Value = 10000
$ws.Cells.Item(110, 13).Value = 9500
$ws.Cells.Item(110, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(110, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(110, 16).Value = 475
$ws.Cells.Item(110, 17).Value = 20
$ws.Cells.Item(110, 18).Value = "Hortaliza"

$ws.Cells.Item(111, 1).Value = 1
$ws.Cells.Item(111, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(111, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(111, 4).Value = 44442
$ws.Cells.Item(111, 4).NumberFormat = $ws.Cells.Item(109, 4).NumberFormat
$ws.Cells.Item(111, 5).Value = 15
$ws.Cells.Item(111, 6).Value = 100112036
$ws.Cells.Item(111, 7).Value = "Caigua"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Segunda"
$ws.Cells.Item(111, 10).Value = 120
$ws.Cells.Item(111, 11).Value = 7000
$ws.Cells.Item(111, 12).Value = 8000
$ws.Cells.Item(111, 13).Value = 7500
$ws.Cells.Item(111, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(111, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(111, 16).Value = 375
$ws.Cells.Item(111, 17).Value = 20
$ws.Cells.Item(111, 18).Value = "Hortaliza"
